# Clientes.xlsx — "Importación y exportación Clientes"
#
# The customer order-history column (D) previously held free-form notes
# ("4 pies", "5 tartaletas", "3 queques", "2 tortaas", "1 torta"). They are
# replaced with a normalized "<n>-<item>," format that matches each row's
# position in the list (row 2 -> item 1, row 3 -> item 2, ...).
#
# Headers (row 1: Nombre / Correo / Número telefónico / Historial de
# pedidos) and the other columns (A, B, C) are untouched — only the shared
# string table is re-packed by the engine once the old strings are no
# longer referenced anywhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "1-torta,"
$ws.Range("D3").Value = "2-torta,"
$ws.Range("D4").Value = "3-queque,"
$ws.Range("D5").Value = "4-pie de limon,"
$ws.Range("D6").Value = "5-tartaleta durazno,"

# Page setup touched (orientation explicitly set to portrait).
$ws.PageSetup.Orientation = 1

# Selection moved from D9 to D10 before the file was saved.
[void]$ws.Range("D10").Select()
